# Apply the perturbation-test edit to the "optimization_parameters" sheet,
# and switch the active/selected sheet from "optimization_diagnostics" to
# "optimization_parameters".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# 1. Row 1 header: remove the duplicated "value" cells in C1:F1 (keep A1/B1).
$ws.Range("C1:F1").ClearContents()

# 2. Rename the "Model" label (row 8, column A) to "production_function".
#    (column B keeps its "Sigmoid" value)
$ws.Range("A8").Value = "production_function"

# 3. Insert a new row directly below it for the "L_curve" parameter.
$ws.Range("A9").EntireRow.Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("B9").NumberFormat = "0.00E+00"

# 4. Remove the old "Deletion" row entirely (it has shifted down to row 17
#    because of the insert above).
$ws.Range("A17").EntireRow.Delete()

# 5. Make this sheet the active / selected sheet, with the last data row
#    (simulation_timepoints, row 17) selected in full.
$ws.Activate()
$ws.Rows.Item(17).Select()
